$d = $word.ActiveDocument

# Locate the paragraph that contains the footnote-reference bullet
# ("Create footnotes with ^[]<footnoteRef 28>") — the new content goes
# right after it, before the "House Price Index" Heading1 paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Create footnotes with*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "could not find the 'Create footnotes with' paragraph"
}

$pRange = $target.Range
# Insert just before the paragraph's own end-of-paragraph mark so the new
# content becomes new sibling paragraphs right after it, instead of
# merging into the following Heading1 paragraph.
$insertAt = $d.Range($pRange.End - 1, $pRange.End - 1)

# The second new paragraph holds a standalone (display) equation. If the
# <m:oMath> is inserted as the paragraph's only content, Word auto-wraps
# it in <m:oMathPara>; the target markup wants a bare <m:oMath> instead,
# so a throwaway empty run is placed immediately before it to suppress
# that auto-wrap, then removed again right after insertion.
$fragment = '<w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1004"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Create math equations with LaTeX syntax:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r></w:r><m:oMath><m:r><m:t>f</m:t></m:r><m:r><m:t>(</m:t></m:r><m:r><m:t>k</m:t></m:r><m:r><m:t>)</m:t></m:r><m:r><m:t>=</m:t></m:r><m:d><m:dPr><m:begChr m:val="("/><m:endChr m:val=")"/><m:grow/></m:dPr><m:e><m:f><m:fPr><m:type m:val="noBar"/></m:fPr><m:num><m:r><m:t>n</m:t></m:r></m:num><m:den><m:r><m:t>k</m:t></m:r></m:den></m:f></m:e></m:d><m:sSup><m:e><m:r><m:t>p</m:t></m:r></m:e><m:sup><m:r><m:t>k</m:t></m:r></m:sup></m:sSup><m:r><m:t>(</m:t></m:r><m:r><m:t>1</m:t></m:r><m:r><m:t>−</m:t></m:r><m:r><m:t>p</m:t></m:r><m:sSup><m:e><m:r><m:t>)</m:t></m:r></m:e><m:sup><m:r><m:t>n</m:t></m:r><m:r><m:t>−</m:t></m:r><m:r><m:t>k</m:t></m:r></m:sup></m:sSup></m:oMath></w:p>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math">' +
    '<w:body>' + $fragment + '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$insertAt.InsertXML($xml)

# Locate the newly-inserted math paragraph (directly follows the new
# "Create math equations..." bullet) and strip the placeholder empty run
# that precedes its <m:oMath>, leaving the bare element as in the target.
$mathPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Create math equations*") {
        $mathPara = $p.Next()
        break
    }
}
if ($mathPara -eq $null) {
    throw "could not find the inserted math paragraph"
}
$placeholder = $d.Range($mathPara.Range.Start, $mathPara.Range.Start + 1)
$placeholder.Delete()

Write-Output "inserted slides bullet + LaTeX oMath example"
